$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 6-9 entirely (cluster now only has 4 control points)
$ws.Range("A6:B9").EntireRow.Delete()

# Update the remaining control-point values in rows 2-5
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 204

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 145

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 104

$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 16
